$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The account-statement table previously listed two workers:
#   - DANNIS DANIEL ARRIETA PEREZ (73552252) with 14 overdue periods
#   - ROBINSON ANTONIO ANTONIO HERNANDEZ ARRIETA (73207993) with 9 periods
# The refreshed export drops the first worker entirely and keeps only
# the second worker's 9 periods, now listed in ascending order
# (2001 .. 2009) with updated "Valor Mora" amounts.
# ------------------------------------------------------------------

# Remove the 14 rows belonging to DANNIS DANIEL ARRIETA PEREZ (rows 16-29).
# Everything below (the remaining worker's rows, the blank spacer rows and
# the signature block) shifts up by 14 rows automatically.
$ws.Range("B16:J29").EntireRow.Delete()

# Rewrite the surviving worker's 9 rows (now B16:J24) in ascending period
# order with the refreshed "Valor Mora" (F) values.
$periods = @(2001, 2002, 2003, 2004, 2005, 2006, 2007, 2008, 2009)
$valorMora = @(1170, 35112, 35112, 35112, 35112, 35112, 35112, 35112, 35112)
$salario = 877803

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = $valorMora[$i]
    $ws.Range("G$r").Value = $salario
}

# Update the summary header figures to match the refreshed table.
$ws.Range("E11").Value = 282066   # VALOR MORA
$ws.Range("C13").Value = 1        # Cant. Trabajadores
$ws.Range("F13").Value = 9        # Cant. Periodos
